$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "Vijayapura (Bijapur)"
$ws.Range("G7").Value = "Bagalkot"
$ws.Range("G16").Value = "Shivamogga (Shimoga)"
$ws.Range("G17").Value = "Bagalkot"
$ws.Range("G21").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G27").Value = "Shivamogga (Shimoga)"
$ws.Range("G29").Value = "Shivamogga (Shimoga)"
$ws.Range("G30").Value = "Davangere"
$ws.Range("G33").Value = "Vijayapura (Bijapur)"
$ws.Range("G34").Value = "Shivamogga (Shimoga)"
$ws.Range("G36").Value = "Shivamogga (Shimoga)"
$ws.Range("G38").Value = "Shivamogga (Shimoga)"
$ws.Range("G39").Value = "Shivamogga (Shimoga)"
$ws.Range("G42").Value = "Bagalkot"
